$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = "armut"
$ws.Range("E5").Value = "beseitigung, verringerung, vermeidung"

$ws.Range("E5").Select()
